$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'iPhone 15 Silicone Case with MagSafe - Winter Blue ​​​​​​​'
$ws.Range("B2").Value = 'B0CHX1652P'
$ws.Range("C2").Value = '''$49.00'
$ws.Range("D2").Value = '''$49.00'
$ws.Range("E2").Value = '4.6 out of 5 stars'
$ws.Range("F2").Value = '''518'
$ws.Range("G2").Value = 'https://www.amazon.com/sspa/click?ie=UTF8&spc=MToxMDg4ODcwMTgxNjM3MDY1OjE3MDk3Mjc5OTU6c3BfYXRmOjMwMDEwNTQ2MTAzMDQwMjo6MDo6&url=%2FApple-iPhone-Silicone-Case-MagSafe%2Fdp%2FB0CHX1652P%2Fref%3Dsr_1_1_sspa%3Fdib%3DeyJ2IjoiMSJ9.qyh_pbIPVBxch9dd0Ix18WCWwT3gv1oE4dFE54aClru3yAF9CQMDr-Z6Rf_g5QahuK0qjdQk-GYxK6AWr2VifgxEnUyRQGCxMvZLu3IzVY_GD6B2fbnG7cJsXjhd7XEW7dbxJVeFdhjfwuU9n6hYmHGUJqCfIA_7zg0qzCSB5OSu6uZ4dAavpoH7yQYhX8UCpooSKlLD7PvccSMygfQu-FvDR8cnVT5cxfHxw6a_TXQ.VIOByNiwmgmOR11EGV5L_gEUHctL0gkY2PvxBb2FmRA%26dib_tag%3Dse%26keywords%3Diphone%2B15%26qid%3D1709727995%26sr%3D8-1-spons%26sp_csd%3Dd2lkZ2V0TmFtZT1zcF9hdGY%26psc%3D1'
$ws.Range("H2").Value = 'https://m.media-amazon.com/images/I/51XqsbKijCL._AC_UY218_.jpg'

# Row 3
$ws.Range("A3").Value = 'iPhone 15 Clear Case with MagSafe ​​​​​​​'
$ws.Range("B3").Value = 'B0CHX1M27P'
$ws.Range("C3").Value = 'Sin Información'
$ws.Range("D3").Value = 'Sin Información'
$ws.Range("E3").Value = '4.3 out of 5 stars'
$ws.Range("F3").Value = '''134'
$ws.Range("G3").Value = 'https://www.amazon.com/sspa/click?ie=UTF8&spc=MToxMDg4ODcwMTgxNjM3MDY1OjE3MDk3Mjc5OTU6c3BfYXRmOjMwMDEwNTQ2MTAzMjkwMjo6MDo6&url=%2FApple-iPhone-Clear-Case-MagSafe%2Fdp%2FB0CHX1M27P%2Fref%3Dsr_1_2_sspa%3Fdib%3DeyJ2IjoiMSJ9.qyh_pbIPVBxch9dd0Ix18WCWwT3gv1oE4dFE54aClru3yAF9CQMDr-Z6Rf_g5QahuK0qjdQk-GYxK6AWr2VifgxEnUyRQGCxMvZLu3IzVY_GD6B2fbnG7cJsXjhd7XEW7dbxJVeFdhjfwuU9n6hYmHGUJqCfIA_7zg0qzCSB5OSu6uZ4dAavpoH7yQYhX8UCpooSKlLD7PvccSMygfQu-FvDR8cnVT5cxfHxw6a_TXQ.VIOByNiwmgmOR11EGV5L_gEUHctL0gkY2PvxBb2FmRA%26dib_tag%3Dse%26keywords%3Diphone%2B15%26qid%3D1709727995%26sr%3D8-2-spons%26sp_csd%3Dd2lkZ2V0TmFtZT1zcF9hdGY%26psc%3D1'
$ws.Range("H3").Value = 'https://m.media-amazon.com/images/I/410mLbFKZRL._AC_UY218_.jpg'

# Row 4
$ws.Range("A4").Value = 'Kryptall 14 Pro Unlocked Cellular Phone, Purple, 128GB Memory Storage Capacity'
$ws.Range("B4").Value = 'B0CJ8B8XJT'
$ws.Range("C4").Value = 'Sin Información'
$ws.Range("D4").Value = 'Sin Información'
$ws.Range("E4").Value = 'Sin Información'
$ws.Range("F4").Value = 'Sin Información'
$ws.Range("G4").Value = 'https://www.amazon.com/Kryptall-Unlocked-Cellular-Storage-Capacity/dp/B0CJ8B8XJT/ref=sr_1_3?dib=eyJ2IjoiMSJ9.qyh_pbIPVBxch9dd0Ix18WCWwT3gv1oE4dFE54aClru3yAF9CQMDr-Z6Rf_g5QahuK0qjdQk-GYxK6AWr2VifgxEnUyRQGCxMvZLu3IzVY_GD6B2fbnG7cJsXjhd7XEW7dbxJVeFdhjfwuU9n6hYmHGUJqCfIA_7zg0qzCSB5OSu6uZ4dAavpoH7yQYhX8UCpooSKlLD7PvccSMygfQu-FvDR8cnVT5cxfHxw6a_TXQ.VIOByNiwmgmOR11EGV5L_gEUHctL0gkY2PvxBb2FmRA&dib_tag=se&keywords=iphone+15&qid=1709727995&sr=8-3'
$ws.Range("H4").Value = 'https://m.media-amazon.com/images/I/61RklF9NgpL._AC_UY218_.jpg'

# Row 5
$ws.Range("A5").Value = 'Galaxy S24+ Plus Cell Phone, 256GB AI Smartphone, Unlocked Android, 50MP Camera, Fastest Processor, Long Battery Life, US Version, 2024, Onyx Black'
$ws.Range("B5").Value = 'B0CMDL3H3V'
$ws.Range("C5").Value = '''$999.99'
$ws.Range("D5").Value = '''$999.99'
$ws.Range("E5").Value = '4.0 out of 5 stars'
$ws.Range("F5").Value = '''338'
$ws.Range("G5").Value = 'https://www.amazon.com/SAMSUNG-Smartphone-Unlocked-Android-Processor/dp/B0CMDL3H3V/ref=sr_1_4?dib=eyJ2IjoiMSJ9.qyh_pbIPVBxch9dd0Ix18WCWwT3gv1oE4dFE54aClru3yAF9CQMDr-Z6Rf_g5QahuK0qjdQk-GYxK6AWr2VifgxEnUyRQGCxMvZLu3IzVY_GD6B2fbnG7cJsXjhd7XEW7dbxJVeFdhjfwuU9n6hYmHGUJqCfIA_7zg0qzCSB5OSu6uZ4dAavpoH7yQYhX8UCpooSKlLD7PvccSMygfQu-FvDR8cnVT5cxfHxw6a_TXQ.VIOByNiwmgmOR11EGV5L_gEUHctL0gkY2PvxBb2FmRA&dib_tag=se&keywords=iphone+15&qid=1709727995&sr=8-4'
$ws.Range("H5").Value = 'https://m.media-amazon.com/images/I/71NngboUC6L._AC_UY218_.jpg'

# Row 6
$ws.Range("A6").Value = 'I15 Pro Max Smartphone Unlocked Cell Phone,Battery 6800mAh 6.8 HD Screen Unlocked Phone,6+256GB Android 13 with 128GB Memory Card,Dual SIM/5G/Fingerprint Lock/Face ID (Purple, 6+256)'
$ws.Range("B6").Value = 'B0CT63GF7X'
$ws.Range("C6").Value = 'Sin Información'
$ws.Range("D6").Value = 'Sin Información'
$ws.Range("E6").Value = '4.7 out of 5 stars'
$ws.Range("F6").Value = '''6'
$ws.Range("G6").Value = 'https://www.amazon.com/Huness-Smartphone-Unlocked-Battery-Fingerprint/dp/B0CT63GF7X/ref=sr_1_5?dib=eyJ2IjoiMSJ9.qyh_pbIPVBxch9dd0Ix18WCWwT3gv1oE4dFE54aClru3yAF9CQMDr-Z6Rf_g5QahuK0qjdQk-GYxK6AWr2VifgxEnUyRQGCxMvZLu3IzVY_GD6B2fbnG7cJsXjhd7XEW7dbxJVeFdhjfwuU9n6hYmHGUJqCfIA_7zg0qzCSB5OSu6uZ4dAavpoH7yQYhX8UCpooSKlLD7PvccSMygfQu-FvDR8cnVT5cxfHxw6a_TXQ.VIOByNiwmgmOR11EGV5L_gEUHctL0gkY2PvxBb2FmRA&dib_tag=se&keywords=iphone+15&qid=1709727995&sr=8-5'
$ws.Range("H6").Value = 'https://m.media-amazon.com/images/I/71UO5Rrpq-L._AC_UY218_.jpg'

# Row 7
$ws.Range("A7").Value = 'Galaxy S23 FE Cell Phone, 256GB, Unlocked Android Smartphone, Long Battery Life, Premium Processor, Tough Gorilla Glass Display, Hi-Res 50MP Camera, US Version, 2023, Cream'
$ws.Range("B7").Value = 'B0CD8YMKWC'
$ws.Range("C7").Value = 'Sin Información'
$ws.Range("D7").Value = 'Sin Información'
$ws.Range("E7").Value = '4.4 out of 5 stars'
$ws.Range("F7").Value = '''540'
$ws.Range("G7").Value = 'https://www.amazon.com/SAMSUNG-Unlocked-Android-Smartphone-Processor/dp/B0CD8YMKWC/ref=sr_1_6?dib=eyJ2IjoiMSJ9.qyh_pbIPVBxch9dd0Ix18WCWwT3gv1oE4dFE54aClru3yAF9CQMDr-Z6Rf_g5QahuK0qjdQk-GYxK6AWr2VifgxEnUyRQGCxMvZLu3IzVY_GD6B2fbnG7cJsXjhd7XEW7dbxJVeFdhjfwuU9n6hYmHGUJqCfIA_7zg0qzCSB5OSu6uZ4dAavpoH7yQYhX8UCpooSKlLD7PvccSMygfQu-FvDR8cnVT5cxfHxw6a_TXQ.VIOByNiwmgmOR11EGV5L_gEUHctL0gkY2PvxBb2FmRA&dib_tag=se&keywords=iphone+15&qid=1709727995&sr=8-6'
$ws.Range("H7").Value = 'https://m.media-amazon.com/images/I/71ukK41npyL._AC_UY218_.jpg'

# Row 8
$ws.Range("A8").Value = 'I15 PROMAX Unlocked Android Phone 2023 Android 13 Cell Phone with Dynamic Island Titanium Design 8GB+512GB Mobile Phones 6.7“ HD Screen 108MP+48MP Camera 6800 mAh Dual SIM Smart Phone (Blue)'
$ws.Range("B8").Value = 'B0CP5K7ZZZ'
$ws.Range("C8").Value = '''$199.99'
$ws.Range("D8").Value = '''$199.99'
$ws.Range("E8").Value = '1.5 out of 5 stars'
$ws.Range("F8").Value = '''11'
$ws.Range("G8").Value = 'https://www.amazon.com/Deeptick-Unlocked-Android-Dynamic-Titanium/dp/B0CP5K7ZZZ/ref=sr_1_7?dib=eyJ2IjoiMSJ9.qyh_pbIPVBxch9dd0Ix18WCWwT3gv1oE4dFE54aClru3yAF9CQMDr-Z6Rf_g5QahuK0qjdQk-GYxK6AWr2VifgxEnUyRQGCxMvZLu3IzVY_GD6B2fbnG7cJsXjhd7XEW7dbxJVeFdhjfwuU9n6hYmHGUJqCfIA_7zg0qzCSB5OSu6uZ4dAavpoH7yQYhX8UCpooSKlLD7PvccSMygfQu-FvDR8cnVT5cxfHxw6a_TXQ.VIOByNiwmgmOR11EGV5L_gEUHctL0gkY2PvxBb2FmRA&dib_tag=se&keywords=iphone+15&qid=1709727995&sr=8-7'
$ws.Range("H8").Value = 'https://m.media-amazon.com/images/I/71v6aolL7qL._AC_UY218_.jpg'

# Row 9
$ws.Range("A9").Value = 'Open, 16GB RAM+512GB, Dual-SIM, Emerald Dusk, US Factory Unlocked Android Smartphone, 4805 mAh Battery, 67W Fast Charging, Hasselblad Camera, 120Hz Fluid Display'
$ws.Range("B9").Value = 'B0CHN8FNW3'
$ws.Range("C9").Value = '''$1,699.99'
$ws.Range("D9").Value = '''$1,699.99'
$ws.Range("E9").Value = '4.0 out of 5 stars'
$ws.Range("F9").Value = '''96'
$ws.Range("G9").Value = 'https://www.amazon.com/OnePlus-Dual-SIM-Unlocked-Smartphone-Hasselblad/dp/B0CHN8FNW3/ref=sr_1_8?dib=eyJ2IjoiMSJ9.qyh_pbIPVBxch9dd0Ix18WCWwT3gv1oE4dFE54aClru3yAF9CQMDr-Z6Rf_g5QahuK0qjdQk-GYxK6AWr2VifgxEnUyRQGCxMvZLu3IzVY_GD6B2fbnG7cJsXjhd7XEW7dbxJVeFdhjfwuU9n6hYmHGUJqCfIA_7zg0qzCSB5OSu6uZ4dAavpoH7yQYhX8UCpooSKlLD7PvccSMygfQu-FvDR8cnVT5cxfHxw6a_TXQ.VIOByNiwmgmOR11EGV5L_gEUHctL0gkY2PvxBb2FmRA&dib_tag=se&keywords=iphone+15&qid=1709727995&sr=8-8'
$ws.Range("H9").Value = 'https://m.media-amazon.com/images/I/61HbEkUjV-L._AC_UY218_.jpg'

# Row 10
$ws.Range("A10").Value = 'A15 Pro Max Cell Phone,8GB+512GB Ultra Memory Unlocked Phone,Android 13.0 Smartphone,6800 mAh Battey,6.82-inch HD Screen,Dual SIM, Dual Standby,108MP Camera, 5G Phone.(Black)'
$ws.Range("B10").Value = 'B0CR2TL446'
$ws.Range("C10").Value = 'Sin Información'
$ws.Range("D10").Value = 'Sin Información'
$ws.Range("E10").Value = '2.5 out of 5 stars'
$ws.Range("F10").Value = '''4'
$ws.Range("G10").Value = 'https://www.amazon.com/WV-LeisureMaster-A15-Smartphone-6-82-inch/dp/B0CR2TL446/ref=sr_1_9?dib=eyJ2IjoiMSJ9.qyh_pbIPVBxch9dd0Ix18WCWwT3gv1oE4dFE54aClru3yAF9CQMDr-Z6Rf_g5QahuK0qjdQk-GYxK6AWr2VifgxEnUyRQGCxMvZLu3IzVY_GD6B2fbnG7cJsXjhd7XEW7dbxJVeFdhjfwuU9n6hYmHGUJqCfIA_7zg0qzCSB5OSu6uZ4dAavpoH7yQYhX8UCpooSKlLD7PvccSMygfQu-FvDR8cnVT5cxfHxw6a_TXQ.VIOByNiwmgmOR11EGV5L_gEUHctL0gkY2PvxBb2FmRA&dib_tag=se&keywords=iphone+15&qid=1709727995&sr=8-9'
$ws.Range("H10").Value = 'https://m.media-amazon.com/images/I/71aGs+ZlfML._AC_UY218_.jpg'

